$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''34.234.98'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +0.39%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.790.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +0.05%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D5').Value = '''226.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.45%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.549'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.39%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.04%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''32.26'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +0.44%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.11%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.0690'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -0.15%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.76%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''2.047.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -0.01%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''11.13'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -3.56%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''1.794.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +0.30%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = '''  +0.45%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''34.193.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +0.29%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.10%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''68.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +0.00%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = '''  +3.54%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''246.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.84%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = '''  +0.30%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '''  +0.06%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''4.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +2.19%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +0.64%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''161.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.61%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''7.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.28%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''16.33'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = '''  +0.83%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  +0.19%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D31').Value = '''0.0520'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.03%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''3.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +2.54%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  +3.88%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  -1.54%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''1.444.20'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +1.94%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '''  +9.29%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.663'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +2.72%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = '''VeChain'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = '''0.0191'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -0.25%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = '''TrustWalletToken'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = '''1.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +1.23%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''82.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +1.85%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''2.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +1.18%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''14.03'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +4.33%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '''  +1.34%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''0.921'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +0.30%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = '''  +2.08%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '''  +0.38%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  +0.32%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''1.947.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -0.04%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''105.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -1.78%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  +0.10%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.0₆0128'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -6.90%  '
$ws.Range('E51').Style = 'Normal'
